$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.068.35"
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("D3").Value = "3.458.93"
$ws.Range("E3").Value = "  -1.64%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.96%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.608"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.18%  "
$ws.Range("D9").Value = "3.457.97"
$ws.Range("E9").Value = "  -1.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.136"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.93"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.19%  "
$ws.Range("E12").Value = "  -3.17%  "
$ws.Range("D13").Value = "4.061.62"
$ws.Range("E13").Value = "  -1.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.54"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.51%  "
$ws.Range("E15").Value = "  -2.91%  "
$ws.Range("D16").Value = "67.060.92"
$ws.Range("E16").Value = "  -1.23%  "
$ws.Range("E17").Value = "  -3.65%  "
$ws.Range("D18").Value = "3.459.50"
$ws.Range("E18").Value = "  -1.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.48%  "
$ws.Range("E20").Value = "  -4.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "379.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.56%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  +0.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("E26").Value = "  -2.23%  "
$ws.Range("E27").Value = "  -2.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.95"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.79%  "
$ws.Range("E29").Value = "  -2.49%  "
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.95"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.18%  "
$ws.Range("E32").Value = "  -3.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.94"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("E34").Value = "  -6.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E36").Value = "  -4.40%  "
$ws.Range("E37").Value = "  -5.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "160.21"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.877"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "27.18"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.81"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.85%  "
$ws.Range("E42").Value = "  -5.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.87%  "
$ws.Range("E44").Value = "  -5.07%  "
$ws.Range("D45").Value = "2.700.03"
$ws.Range("E45").Value = "  -6.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0699"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "41.05"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.10%  "
$ws.Range("E49").Value = "  -3.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "322.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.03"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.26%  "
